$wb = $excel.ActiveWorkbook

# --- Update status text "Ready for handoff" -> "In Translation" ---
# Sheet "Overview": columns E (zh-cn) and F (de-de), rows 2-3
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Range("E3").Value = "In Translation"
$wsOverview.Range("F3").Value = "In Translation"

# Sheet "zh-cn": Status column C, rows 2-3
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Range("C3").Value = "In Translation"

# Sheet "de-de": Status column C, rows 2-3
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Range("C3").Value = "In Translation"

# --- Narrow the "Status" columns (E/F on Overview, C on zh-cn / de-de) ---
$wsOverview.Columns.Item(5).ColumnWidth = 12.45
$wsOverview.Columns.Item(6).ColumnWidth = 12.45
$wsZhCn.Columns.Item(3).ColumnWidth = 12.45
$wsDeDe.Columns.Item(3).ColumnWidth = 12.45
